# Auto-generated Excel COM-interop script to apply cryptos.xlsx data refresh
# (crypto price/volume table update, commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.983.35'
$ws.Range('D3').Value = '1.642.11'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5062'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2580'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06366'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.89'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07732'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.299'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').Value = '1.645.32'
$ws.Range('E13').Value = '  +0.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5472'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.41%  '
$ws.Range('D15').Value = '0.0₅7763'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.28'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('D17').Value = '26.012.49'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.473'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '196.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.983'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.152'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.003'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.893'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.53'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1265'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.68%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.876'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.64'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.240'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04906'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.271'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.208'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.552'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.378'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9192'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.87%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '1.138.95'
$ws.Range('E36').Value = '  +0.87%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.568'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.64%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5550'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.002'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.608'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.49%  '
$ws.Range('E42').Value = '  -1.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '98.80'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.21%  '
$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D44').Value = '0.0₈120'
$ws.Range('E44').Value = '  -8.38%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.778.32'
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4527'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.37'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.37%  '
$ws.Range('E49').Value = '  +2.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.584'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.53%  '
$ws.Range('E51').Value = '  -0.24%  '
